$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: rename dimension-style labels to measure-style, and repurpose H2
$ws.Range("E2").Value = "iaest-measure:codcom"
$ws.Range("G2").Value = "iaest-measure:sector-descripcion"
$ws.Range("H2").Value = "iaest-measure:direccion-provincial-nombre"

# Row 3: "dim" -> "medida" for the columns that changed from dimension to measure
$ws.Range("E3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4: datatype updates
$ws.Range("E4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"

# Row 5 (mapping file references) no longer needed - remove entirely
$ws.Range("A5:I5").EntireRow.Delete()
